$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 9.4188686742162915
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 10.054387545355425
$ws.Range("E2").Value = 8.2091894300565773

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 7.5300292502696546
$ws.Range("D3").Value = 6.039703408895261
$ws.Range("E3").Value = 8.5190974215319812

# Update selection to reflect the narrower range used
$ws.Range("B1:E3").Select()
